$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before O (current syst_u) -> new column becomes O
# ("scaled st"), and the old O (syst_u) shifts right to P.
$ws.Columns.Item(15).Insert()

# Insert a new column after the (now shifted) syst_u column P -> new column
# becomes Q ("scaled sy").
$ws.Columns.Item(17).Insert()

# Headers
$ws.Range("O1").Value = "scaled st"
$ws.Range("Q1").Value = "scaled sy"

# Formulas for the scaled columns. Row 2 is a standalone formula, while
# rows 3:66 and 67:82 are filled as ranges (two separate data blocks),
# matching the original M-column (value) shared-formula grouping.
$ws.Range("O2").Formula = "=N2/100"
$ws.Range("O3:O66").Formula = "=N3/100"
$ws.Range("O67:O82").Formula = "=N67/100"

$ws.Range("Q2").Formula = "=P2/100"
$ws.Range("Q3:Q66").Formula = "=P3/100"
$ws.Range("Q67:Q82").Formula = "=P67/100"

# View adjustments matching the diff (scroll + selection)
$ws.Application.ActiveWindow.ScrollRow = 64
$ws.Range("Q2:Q82").Select()
